$wb = $excel.ActiveWorkbook

# ALC (sheet 1), row 17
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 1346.25  # H17: 1356 -> 1346.25
$ws.Cells.Item(17, 10).Value = 1365.1333  # J17: 1376.9286 -> 1365.1333
$ws.Cells.Item(17, 12).Value = 4095.3999  # L17: 4130.7858 -> 4095.3999
$ws.Cells.Item(17, 14).Value = -4431.3999  # N17: -4466.7858 -> -4431.3999

# ALC (sheet 1), row 92
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(92, 8).Value = 537.6667  # H92: 610.7692 -> 537.6667
$ws.Cells.Item(92, 9).Value = 238.75  # I92: 274 -> 238.75
$ws.Cells.Item(92, 11).Value = 238.75  # K92: 274 -> 238.75
$ws.Cells.Item(92, 13).Value = 1009.25  # M92: 974 -> 1009.25

# ALC (sheet 1), row 111
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(111, 8).Value = 807.46155  # H111: 1077.1 -> 807.46155
$ws.Cells.Item(111, 9).Value = 964.1  # I111: 1077.1 -> 964.1
$ws.Cells.Item(111, 10).Value = 285.33334  # J111: 0 -> 285.33334
$ws.Cells.Item(111, 11).Value = 2892.3  # K111: 3231.3 -> 2892.3
$ws.Cells.Item(111, 12).Value = 856.0000200000001  # L111: 0 -> 856.0000200000001
$ws.Cells.Item(111, 13).Value = 174.6999999999998  # M111: -164.2999999999997 -> 174.6999999999998
$ws.Cells.Item(111, 14).Value = -6990.00002  # N111: None -> -6990.00002

# ALC (sheet 1), row 116
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(116, 8).Value = 2999.8572  # H116: 2874.875 -> 2999.8572
$ws.Cells.Item(116, 10).Value = 5999  # J116: 3999.5 -> 5999
$ws.Cells.Item(116, 12).Value = 5999  # L116: 3999.5 -> 5999
$ws.Cells.Item(116, 14).Value = -12883  # N116: -10883.5 -> -12883

# ALC (sheet 1), row 138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 4284.229  # H138: 4364.1226 -> 4284.229
$ws.Cells.Item(138, 10).Value = 4924.636  # J138: 5067 -> 4924.636
$ws.Cells.Item(138, 12).Value = 14773.908  # L138: 15201 -> 14773.908
$ws.Cells.Item(138, 14).Value = -25053.908  # N138: -25481 -> -25053.908

# ARM (sheet 2), row 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 1160.875  # H2: 1221.1818 -> 1160.875
$ws.Cells.Item(2, 9).Value = 626.8095  # I2: 640.4211 -> 626.8095
$ws.Cells.Item(2, 11).Value = 626.8095  # K2: 640.4211 -> 626.8095
$ws.Cells.Item(2, 13).Value = -513.8095  # M2: -527.4211 -> -513.8095

# ARM (sheet 2), row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 1499.2858  # H61: 1384.8182 -> 1499.2858
$ws.Cells.Item(61, 9).Value = 1499.2858  # I61: 1433.3 -> 1499.2858
$ws.Cells.Item(61, 10).Value = 0  # J61: 900 -> 0
$ws.Cells.Item(61, 11).Value = 1499.2858  # K61: 1433.3 -> 1499.2858
$ws.Cells.Item(61, 12).Value = 0  # L61: 900 -> 0
$ws.Cells.Item(61, 13).Value = -1287.2858  # M61: -1221.3 -> -1287.2858
$ws.Cells.Item(61, 14).ClearContents()  # N61: -1324 -> (removed)

# ARM (sheet 2), row 74
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 1214.4615  # H74: 1221.4615 -> 1214.4615
$ws.Cells.Item(74, 9).Value = 1013.8  # I74: 1022.9 -> 1013.8
$ws.Cells.Item(74, 11).Value = 1013.8  # K74: 1022.9 -> 1013.8
$ws.Cells.Item(74, 13).Value = -139.8  # M74: -148.9 -> -139.8

# ARM (sheet 2), row 77
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value = 1214.4615  # H77: 1221.4615 -> 1214.4615
$ws.Cells.Item(77, 9).Value = 1013.8  # I77: 1022.9 -> 1013.8
$ws.Cells.Item(77, 11).Value = 5069  # K77: 5114.5 -> 5069
$ws.Cells.Item(77, 13).Value = -701  # M77: -746.5 -> -701

# ARM (sheet 2), row 97
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(97, 8).Value = 518.5714  # H97: 405 -> 518.5714
$ws.Cells.Item(97, 10).Value = 1200  # J97: 0 -> 1200
$ws.Cells.Item(97, 12).Value = 1200  # L97: 0 -> 1200
$ws.Cells.Item(97, 14).Value = -2192  # N97: None -> -2192

# ARM (sheet 2), row 116
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(116, 8).Value = 1160.875  # H116: 1221.1818 -> 1160.875
$ws.Cells.Item(116, 9).Value = 626.8095  # I116: 640.4211 -> 626.8095
$ws.Cells.Item(116, 11).Value = 626.8095  # K116: 640.4211 -> 626.8095
$ws.Cells.Item(116, 13).Value = 1667.1905  # M116: 1653.5789 -> 1667.1905

# ARM (sheet 2), row 122
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value = 2597.875  # H122: 2597.8125 -> 2597.875
$ws.Cells.Item(122, 9).Value = 2312.8462  # I122: 2361.7856 -> 2312.8462
$ws.Cells.Item(122, 10).Value = 3833  # J122: 4250 -> 3833
$ws.Cells.Item(122, 11).Value = 6938.5386  # K122: 7085.3568 -> 6938.5386
$ws.Cells.Item(122, 12).Value = 11499  # L122: 12750 -> 11499
$ws.Cells.Item(122, 13).Value = -4488.5386  # M122: -4635.3568 -> -4488.5386
$ws.Cells.Item(122, 14).Value = -16399  # N122: -17650 -> -16399

# ARM (sheet 2), row 132
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 904  # H132: 1056.3334 -> 904
$ws.Cells.Item(132, 9).Value = 622.6667  # I132: 646 -> 622.6667
$ws.Cells.Item(132, 11).Value = 1868.0001  # K132: 1938 -> 1868.0001
$ws.Cells.Item(132, 13).Value = 661.9999  # M132: 592 -> 661.9999

# ARM (sheet 2), row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 1499.2858  # H136: 1384.8182 -> 1499.2858
$ws.Cells.Item(136, 9).Value = 1499.2858  # I136: 1433.3 -> 1499.2858
$ws.Cells.Item(136, 10).Value = 0  # J136: 900 -> 0
$ws.Cells.Item(136, 11).Value = 4497.857400000001  # K136: 4299.9 -> 4497.857400000001
$ws.Cells.Item(136, 12).Value = 0  # L136: 2700 -> 0
$ws.Cells.Item(136, 13).Value = -1947.857400000001  # M136: -1749.9 -> -1947.857400000001
$ws.Cells.Item(136, 14).ClearContents()  # N136: -7800 -> (removed)

# BSM (sheet 3), row 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 1160.875  # H3: 1221.1818 -> 1160.875
$ws.Cells.Item(3, 9).Value = 626.8095  # I3: 640.4211 -> 626.8095
$ws.Cells.Item(3, 11).Value = 626.8095  # K3: 640.4211 -> 626.8095
$ws.Cells.Item(3, 13).Value = -512.8095  # M3: -526.4211 -> -512.8095

# BSM (sheet 3), row 99
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(99, 8).Value = 2853.1333  # H99: 2853.2 -> 2853.1333
$ws.Cells.Item(99, 9).Value = 3066  # I99: 3299.5 -> 3066
$ws.Cells.Item(99, 10).Value = 2799.9167  # J99: 2784.5386 -> 2799.9167
$ws.Cells.Item(99, 11).Value = 3066  # K99: 3299.5 -> 3066
$ws.Cells.Item(99, 12).Value = 2799.9167  # L99: 2784.5386 -> 2799.9167
$ws.Cells.Item(99, 13).Value = -1568  # M99: -1801.5 -> -1568
$ws.Cells.Item(99, 14).Value = -5795.9167  # N99: -5780.5386 -> -5795.9167

# BSM (sheet 3), row 134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 3570.5715  # H134: 3749.3076 -> 3570.5715
$ws.Cells.Item(134, 9).Value = 2937.25  # I134: 3166.6667 -> 2937.25
$ws.Cells.Item(134, 10).Value = 3823.9  # J134: 3924.1 -> 3823.9
$ws.Cells.Item(134, 11).Value = 8811.75  # K134: 9500.000100000001 -> 8811.75
$ws.Cells.Item(134, 12).Value = 11471.7  # L134: 11772.3 -> 11471.7
$ws.Cells.Item(134, 13).Value = -6276.75  # M134: -6965.000100000001 -> -6276.75
$ws.Cells.Item(134, 14).Value = -16541.7  # N134: -16842.3 -> -16541.7

# CRP (sheet 4), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2802.182  # H31: 2881.879 -> 2802.182
$ws.Cells.Item(31, 9).Value = 2292.7144  # I31: 2386.6428 -> 2292.7144
$ws.Cells.Item(31, 11).Value = 2292.7144  # K31: 2386.6428 -> 2292.7144
$ws.Cells.Item(31, 13).Value = -1997.7144  # M31: -2091.6428 -> -1997.7144

# CRP (sheet 4), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 2802.182  # H34: 2881.879 -> 2802.182
$ws.Cells.Item(34, 9).Value = 2292.7144  # I34: 2386.6428 -> 2292.7144
$ws.Cells.Item(34, 11).Value = 2292.7144  # K34: 2386.6428 -> 2292.7144
$ws.Cells.Item(34, 13).Value = -2090.7144  # M34: -2184.6428 -> -2090.7144

# CRP (sheet 4), row 62
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(62, 8).Value = 35588.668  # H62: 47040.555 -> 35588.668
$ws.Cells.Item(62, 9).Value = 2210.8333  # I62: 2816.5 -> 2210.8333
$ws.Cells.Item(62, 10).Value = 68966.5  # J62: 82419.8 -> 68966.5
$ws.Cells.Item(62, 11).Value = 2210.8333  # K62: 2816.5 -> 2210.8333
$ws.Cells.Item(62, 12).Value = 68966.5  # L62: 82419.8 -> 68966.5
$ws.Cells.Item(62, 13).Value = -1586.8333  # M62: -2192.5 -> -1586.8333
$ws.Cells.Item(62, 14).Value = -70214.5  # N62: -83667.8 -> -70214.5

# CRP (sheet 4), row 65
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(65, 8).Value = 35588.668  # H65: 47040.555 -> 35588.668
$ws.Cells.Item(65, 9).Value = 2210.8333  # I65: 2816.5 -> 2210.8333
$ws.Cells.Item(65, 10).Value = 68966.5  # J65: 82419.8 -> 68966.5
$ws.Cells.Item(65, 11).Value = 11054.1665  # K65: 14082.5 -> 11054.1665
$ws.Cells.Item(65, 12).Value = 344832.5  # L65: 412099 -> 344832.5
$ws.Cells.Item(65, 13).Value = -7934.166499999999  # M65: -10962.5 -> -7934.166499999999
$ws.Cells.Item(65, 14).Value = -351072.5  # N65: -418339 -> -351072.5

# CRP (sheet 4), row 99
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(99, 8).Value = 13502.739  # H99: 14467.19 -> 13502.739
$ws.Cells.Item(99, 9).Value = 10615.889  # I99: 11566.375 -> 10615.889
$ws.Cells.Item(99, 10).Value = 15358.571  # J99: 16252.308 -> 15358.571
$ws.Cells.Item(99, 11).Value = 10615.889  # K99: 11566.375 -> 10615.889
$ws.Cells.Item(99, 12).Value = 15358.571  # L99: 16252.308 -> 15358.571
$ws.Cells.Item(99, 13).Value = -9117.888999999999  # M99: -10068.375 -> -9117.888999999999
$ws.Cells.Item(99, 14).Value = -18354.571  # N99: -19248.308 -> -18354.571

# CRP (sheet 4), row 126
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(126, 8).Value = 13502.739  # H126: 14467.19 -> 13502.739
$ws.Cells.Item(126, 9).Value = 10615.889  # I126: 11566.375 -> 10615.889
$ws.Cells.Item(126, 10).Value = 15358.571  # J126: 16252.308 -> 15358.571
$ws.Cells.Item(126, 11).Value = 31847.667  # K126: 34699.125 -> 31847.667
$ws.Cells.Item(126, 12).Value = 46075.713  # L126: 48756.924 -> 46075.713
$ws.Cells.Item(126, 13).Value = -29377.667  # M126: -32229.125 -> -29377.667
$ws.Cells.Item(126, 14).Value = -51015.713  # N126: -53696.924 -> -51015.713

# CRP (sheet 4), row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 1773.825  # H134: 1868.6945 -> 1773.825
$ws.Cells.Item(134, 9).Value = 1595.8235  # I134: 1685.9333 -> 1595.8235
$ws.Cells.Item(134, 11).Value = 4787.470499999999  # K134: 5057.7999 -> 4787.470499999999
$ws.Cells.Item(134, 13).Value = -2252.470499999999  # M134: -2522.7999 -> -2252.470499999999

# CUL (sheet 5), row 75
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(75, 8).Value = 2350.75  # H75: 2584.5 -> 2350.75
$ws.Cells.Item(75, 9).Value = 1999  # I75: 2498.5 -> 1999
$ws.Cells.Item(75, 10).Value = 2561.8  # J75: 2627.5 -> 2561.8
$ws.Cells.Item(75, 11).Value = 5997  # K75: 7495.5 -> 5997
$ws.Cells.Item(75, 12).Value = 7685.400000000001  # L75: 7882.5 -> 7685.400000000001
$ws.Cells.Item(75, 13).Value = -4999  # M75: -6497.5 -> -4999
$ws.Cells.Item(75, 14).Value = -9681.400000000001  # N75: -9878.5 -> -9681.400000000001

# CUL (sheet 5), row 78
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(78, 8).Value = 2350.75  # H78: 2584.5 -> 2350.75
$ws.Cells.Item(78, 9).Value = 1999  # I78: 2498.5 -> 1999
$ws.Cells.Item(78, 10).Value = 2561.8  # J78: 2627.5 -> 2561.8
$ws.Cells.Item(78, 11).Value = 17991  # K78: 22486.5 -> 17991
$ws.Cells.Item(78, 12).Value = 23056.2  # L78: 23647.5 -> 23056.2
$ws.Cells.Item(78, 13).Value = -12999  # M78: -17494.5 -> -12999
$ws.Cells.Item(78, 14).Value = -33040.2  # N78: -33631.5 -> -33040.2

# CUL (sheet 5), row 103
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(103, 8).Value = 167.33333  # H103: 179.4 -> 167.33333
$ws.Cells.Item(103, 9).Value = 137  # I103: 167.25 -> 137
$ws.Cells.Item(103, 11).Value = 411  # K103: 501.75 -> 411
$ws.Cells.Item(103, 13).Value = 468  # M103: 377.25 -> 468

# CUL (sheet 5), row 118
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(118, 8).Value = 904.75  # H118: 1162.2 -> 904.75
$ws.Cells.Item(118, 9).Value = 873  # I118: 897.25 -> 873
$ws.Cells.Item(118, 10).Value = 1000  # J118: 2222 -> 1000
$ws.Cells.Item(118, 11).Value = 2619  # K118: 2691.75 -> 2619
$ws.Cells.Item(118, 12).Value = 3000  # L118: 6666 -> 3000
$ws.Cells.Item(118, 13).Value = -1376  # M118: -1448.75 -> -1376
$ws.Cells.Item(118, 14).Value = -5486  # N118: -9152 -> -5486

# GSM (sheet 6), row 14
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(14, 8).Value = 1673333  # H14: 0 -> 1673333
$ws.Cells.Item(14, 9).Value = 1673333  # I14: 0 -> 1673333
$ws.Cells.Item(14, 11).Value = 1673333  # K14: 0 -> 1673333
$ws.Cells.Item(14, 13).Value = -1673165  # M14: None -> -1673165

# GSM (sheet 6), row 20
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(20, 8).Value = 11199.8  # H20: 10062.375 -> 11199.8
$ws.Cells.Item(20, 10).Value = 12500  # J20: 10642.857 -> 12500
$ws.Cells.Item(20, 12).Value = 12500  # L20: 10642.857 -> 12500
$ws.Cells.Item(20, 14).Value = -12990  # N20: -11132.857 -> -12990

# GSM (sheet 6), row 107
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(107, 8).Value = 839.8  # H107: 1033.3334 -> 839.8
$ws.Cells.Item(107, 9).Value = 633  # I107: 800 -> 633
$ws.Cells.Item(107, 11).Value = 633  # K107: 800 -> 633
$ws.Cells.Item(107, 13).Value = 1287  # M107: 1120 -> 1287

# GSM (sheet 6), row 126
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(126, 8).Value = 3342.7144  # H126: 2897.5 -> 3342.7144

# GSM (sheet 6), row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 2361.96  # H132: 2412.7917 -> 2361.96
$ws.Cells.Item(132, 9).Value = 2288.5881  # I132: 2360.25 -> 2288.5881
$ws.Cells.Item(132, 11).Value = 6865.7643  # K132: 7080.75 -> 6865.7643
$ws.Cells.Item(132, 13).Value = -4335.7643  # M132: -4550.75 -> -4335.7643

# LTW (sheet 7), row 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 4222.5  # H7: 4922.5 -> 4222.5
$ws.Cells.Item(7, 9).Value = 4000  # I7: 4922.5 -> 4000
$ws.Cells.Item(7, 10).Value = 4445  # J7: 0 -> 4445
$ws.Cells.Item(7, 11).Value = 4000  # K7: 4922.5 -> 4000
$ws.Cells.Item(7, 12).Value = 4445  # L7: 0 -> 4445
$ws.Cells.Item(7, 13).Value = -3888  # M7: -4810.5 -> -3888
$ws.Cells.Item(7, 14).Value = -4669  # N7: None -> -4669

# LTW (sheet 7), row 61
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61, 8).Value = 3383.2856  # H61: 3466.6155 -> 3383.2856
$ws.Cells.Item(61, 9).Value = 3258.923  # I61: 3338.8333 -> 3258.923
$ws.Cells.Item(61, 11).Value = 3258.923  # K61: 3338.8333 -> 3258.923
$ws.Cells.Item(61, 13).Value = -3056.923  # M61: -3136.8333 -> -3056.923

# LTW (sheet 7), row 113
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(113, 8).Value = 3383.2856  # H113: 3466.6155 -> 3383.2856
$ws.Cells.Item(113, 9).Value = 3258.923  # I113: 3338.8333 -> 3258.923
$ws.Cells.Item(113, 11).Value = 3258.923  # K113: 3338.8333 -> 3258.923
$ws.Cells.Item(113, 13).Value = -1088.923  # M113: -1168.8333 -> -1088.923

# LTW (sheet 7), row 126
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(126, 8).Value = 4222.5  # H126: 4922.5 -> 4222.5
$ws.Cells.Item(126, 9).Value = 4000  # I126: 4922.5 -> 4000
$ws.Cells.Item(126, 10).Value = 4445  # J126: 0 -> 4445
$ws.Cells.Item(126, 11).Value = 12000  # K126: 14767.5 -> 12000
$ws.Cells.Item(126, 12).Value = 13335  # L126: 0 -> 13335
$ws.Cells.Item(126, 13).Value = -9530  # M126: -12297.5 -> -9530
$ws.Cells.Item(126, 14).Value = -18275  # N126: None -> -18275

# WVR (sheet 8), row 126
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(126, 8).Value = 5718.25  # H126: 4634.5386 -> 5718.25
$ws.Cells.Item(126, 9).Value = 4349.8  # I126: 3625.2 -> 4349.8
$ws.Cells.Item(126, 11).Value = 13049.4  # K126: 10875.6 -> 13049.4
$ws.Cells.Item(126, 13).Value = -10579.4  # M126: -8405.599999999999 -> -10579.4

# WVR (sheet 8), row 132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 6655.3125  # H132: 7181.0713 -> 6655.3125
$ws.Cells.Item(132, 9).Value = 3909.7778  # I132: 4176.857 -> 3909.7778
$ws.Cells.Item(132, 11).Value = 11729.3334  # K132: 12530.571 -> 11729.3334
$ws.Cells.Item(132, 13).Value = -9199.3334  # M132: -10000.571 -> -9199.3334
